$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (WP3932)
$ws.Range("E2").Value = 1.497902135034229
$ws.Range("F2").Value = 0.04353784861795712
$ws.Range("G2").Value = 0.31238095238095237
$ws.Range("H2").Value = 0.30142021720969087

# Row 3 (WP2877)
$ws.Range("E3").Value = -1.37487132613748
$ws.Range("F3").Value = 0.06476190476190476
$ws.Range("G3").Value = 0.31238095238095237
$ws.Range("H3").Value = 0.30142021720969087

# Row 4 (WP1449)
$ws.Range("E4").Value = -1.3553100576828794
$ws.Range("F4").Value = 0.07809523809523809
$ws.Range("G4").Value = 0.31238095238095237
$ws.Range("H4").Value = 0.30142021720969087
